$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reordered roster data (player, position, team) for rows 2-19
$data = @(
    @("Jalen Green",      "PG,SG",    "Houston Rockets"),
    @("Buddy Hield",      "SG,SF",    "Golden State Warriors"),
    @("Grayson Allen",    "PG,SG,SF", "Phoenix Suns"),
    @("Caris LeVert",     "SG,SF",    "Cleveland Cavaliers"),
    @("Jaylen Brown",     "SG,SF",    "Boston Celtics"),
    @("Rudy Gobert",      "C",        "Minnesota Timberwolves"),
    @("Pascal Siakam",    "SF,PF,C",  "Indiana Pacers"),
    @("Jakob Poeltl",     "C",        "Toronto Raptors"),
    @("Nikola Jokic",     "C",        "Denver Nuggets"),
    @("Paolo Banchero",   "SF,PF",    "Orlando Magic"),
    @("Khris Middleton",  "SF",       "Milwaukee Bucks"),
    @("Chris Paul",       "PG",       "San Antonio Spurs"),
    @("Dejounte Murray",  "PG,SG",    "New Orleans Pelicans"),
    @("Russell Westbrook","PG,SG",    "Denver Nuggets"),
    @("Keon Johnson",     "PG,SG",    "Brooklyn Nets"),
    @("Chet Holmgren",    "PF,C",     "Oklahoma City Thunder"),
    @("Jalen Suggs",      "PG,SG",    "Orlando Magic"),
    @("Deni Avdija",      "SF,PF",    "Portland Trail Blazers")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
